$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.527.34"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "2.487.21"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.63"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.18"
$ws.Range("E6").Value = "  -1.20%  "
$ws.Range("E7").Value = "  -1.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.59"
$ws.Range("E10").Value = "  -2.72%  "
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("E12").Value = "  +2.00%  "
$ws.Range("D13").Value = "2.871.49"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.83"
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("D15").Value = "2.538.63"
$ws.Range("E15").Value = "  +2.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.42"
$ws.Range("E16").Value = "  +5.95%  "
$ws.Range("E17").Value = "  -4.51%  "
$ws.Range("D18").Value = "41.589.97"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.32"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").Value = "0.0₃0927"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.57"
$ws.Range("E21").Value = "  +4.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.20"
$ws.Range("E22").Value = "  -2.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.59"
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("E24").Value = "  -2.96%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.90"
$ws.Range("E26").Value = "  -1.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.90"
$ws.Range("E27").Value = "  +1.52%  "
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.64"
$ws.Range("E29").Value = "  -0.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.24"
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "154.96"
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.42"
$ws.Range("E32").Value = "  -3.21%  "
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.19"
$ws.Range("E34").Value = "  +6.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0755"
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.45"
$ws.Range("E36").Value = "  -4.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.95"
$ws.Range("E37").Value = "  -2.27%  "
$ws.Range("E38").Value = "  -3.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.105"
$ws.Range("E39").Value = "  +1.33%  "
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.10"
$ws.Range("E41").Value = "  -3.47%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.16"
$ws.Range("E43").Value = "  -4.65%  "
$ws.Range("D44").Value = "1.963.94"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("E46").Value = "  -3.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.84"
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("D48").Value = "2.727.42"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "96.27"
$ws.Range("E49").Value = "  -1.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.30"
$ws.Range("E50").Value = "  -3.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.19"
$ws.Range("E51").Value = "  -4.17%  "
